$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.969.19"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.649.81"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3913"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  +3.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.371"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08496"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.254"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.983"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001319"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "1.650.75"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06988"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.948"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "23.955.51"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.474"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.302"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.908"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.488"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "1.831.68"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.049"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.95%  "
$ws.Range("E35").Value = "  +4.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08152"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2726"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09190"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7577"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.424"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6976"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.496"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.099"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08312"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.404"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.39%  "
